# Update cryptos list cell values per upstream data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.862.00"
$ws.Range("D3").Value = "2.533.80"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.21"
$ws.Range("E5").Value = "  +1.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.03"
$ws.Range("E6").Value = "  +2.19%  "
$ws.Range("E7").Value = "  -0.78%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("E10").Value = "  -1.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0816"
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.54"
$ws.Range("E12").Value = "  -2.04%  "
$ws.Range("E13").Value = "  -3.61%  "
$ws.Range("D14").Value = "2.925.57"
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").Value = "2.537.48"
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("E16").Value = "  -3.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.852"
$ws.Range("D18").Value = "42.933.25"
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.93"
$ws.Range("E19").Value = "  +3.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.63"
$ws.Range("E20").Value = "  -3.19%  "
$ws.Range("D21").Value = "0.0₃0967"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.71"
$ws.Range("E22").Value = "  -1.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "252.90"
$ws.Range("E23").Value = "  -0.69%  "
$ws.Range("E24").Value = "  +1.07%  "
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.42"
$ws.Range("E26").Value = "  -3.47%  "
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.41"
$ws.Range("E28").Value = "  +1.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "41.34"
$ws.Range("E29").Value = "  +3.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.49"
$ws.Range("E30").Value = "  +4.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.90"
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.22"
$ws.Range("E32").Value = "  +1.02%  "
$ws.Range("E33").Value = "  +0.61%  "
$ws.Range("B34").Value = "Celestia"
$ws.Range("C34").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.39"
$ws.Range("E34").Value = "  -1.66%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.36"
$ws.Range("E35").Value = "  -1.53%  "
$ws.Range("E36").Value = "  +3.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0793"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("E38").Value = "  +0.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.45"
$ws.Range("E39").Value = "  +13.26%  "
$ws.Range("E40").Value = "  -0.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "21.78"
$ws.Range("E41").Value = "  -12.05%  "
$ws.Range("E42").Value = "  +0.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.82"
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("E45").Value = "  -3.26%  "
$ws.Range("D46").Value = "2.008.18"
$ws.Range("E46").Value = "  -1.56%  "
$ws.Range("E47").Value = "  +2.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "84.27"
$ws.Range("E48").Value = "  -1.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "106.19"
$ws.Range("E49").Value = "  +3.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.18"
$ws.Range("E50").Value = "  -0.81%  "
$ws.Range("D51").Value = "2.779.54"
$ws.Range("E51").Value = "  -0.16%  "
